$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9846.0
$ws.Range("C2").Value = 12947265.0
$ws.Range("D2").Value = 15068.0
$ws.Range("E2").Value = 4260.62
$ws.Range("F2").Value = 190.0
$ws.Range("G2").Value = 103.0
$ws.Range("H2").Value = 41.3652427184466
$ws.Range("I2").Value = 22.424315789473685
